# "Fix de imagenes de los paretos"
# The results table on Hoja1 gains a leading "Metodo" column identifying
# each row (SMARTER, Fuzzy, TOPSIS, GRA, CODAS, MABAC, VIKOR, PROMETHEE II)
# and the remaining headers are renamed to Rx/Ry/CL/Entropia/SSIM.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new column A, shifting the existing A:E data to B:F.
$ws.Columns.Item(1).Insert(-4161)  # xlShiftToRight

# New column A: method header + the 8 method names.
$ws.Range("A1").Value = "Metodo"
$ws.Range("A2").Value = "SMARTER"
$ws.Range("A3").Value = "Fuzzy"
$ws.Range("A4").Value = "TOPSIS"
$ws.Range("A5").Value = "GRA"
$ws.Range("A6").Value = "CODAS"
$ws.Range("A7").Value = "MABAC"
$ws.Range("A8").Value = "VIKOR"
$ws.Range("A9").Value = "PROMETHEE II"

# Former headers (Var1_1..Var1_5, now in B1:F1) get renamed.
$ws.Range("B1").Value = "Rx"
$ws.Range("C1").Value = "Ry"
$ws.Range("D1").Value = "CL"
$ws.Range("E1").Value = "Entropia"
$ws.Range("F1").Value = "SSIM"

# Re-fit the three new columns to their (now narrower/wider) content.
$ws.Columns.Item(1).ColumnWidth = 12.5
$ws.Columns.Item(2).ColumnWidth = 2.3333333333
$ws.Columns.Item(3).ColumnWidth = 3.1666666667

$wb.Save()
